$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update CID values on the existing enemy rows (battle support enemy logic) ---
$ws.Range("B2").Value = 100
$ws.Range("B3").Value = 101
$ws.Range("B4").Value = 102

# --- Prepare rows 5, 6 and 7 so that columns A/B pick up the same cell
# style (s="1") that rows 2-4 already carry. Inserting a row copies the
# formatting down from the row above (A4/B4 -> new row), and the stray
# row that gets pushed down by the insert is removed right afterwards,
# leaving the row count unchanged but the new row properly styled. ---
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Delete()

$ws.Rows.Item(6).Insert()
$ws.Rows.Item(7).Delete()

$ws.Rows.Item(7).Insert()
$ws.Rows.Item(8).Delete()

# --- New enemy definitions (IDs 3, 4, 5 / CIDs 103, 104, 105) ---
$rowsData = @{
    5 = @(3, 103, 1, 0, 0, 0, 0, 1, 0, 1, 0, 2, 0, 2, 0)
    6 = @(4, 104, 2, 1, 0, 0, 0, 2, 0, 1, 0, 2, 0, 2, 0)
    7 = @(5, 105, 0, 2, 0, 0, 0, 0, 0, 1, 0, 1, 0, 2, 0)
}

foreach ($r in 5..7) {
    $values = $rowsData[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}

# --- Bump the sheet's cached outline-level-row high-water mark to 6,
# matching the source workbook's bookkeeping, without leaving a visible
# outlineLevel on any real data row. ---
$ws.Rows.Item(9).OutlineLevel = 6
$ws.Rows.Item(9).Delete()

# --- Move the active selection to reflect the last edited cell ---
$ws.Range("D6").Select()
